$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add "X" marker in column E (Hoàn thành) for rows 2-6
$ws.Range("E2:E6").Value = "X"
